$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.012.24"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.65"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.61"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5130"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06378"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.81"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07813"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.290"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.651.08"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5473"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.54"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7754"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.074.26"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.006"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.06"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.449"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.978"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.084"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.899"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.39"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1223"
$ws.Range("E26").Value = "  +7.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.877"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.242"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04885"
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.283"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.207"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.542"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.384"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9142"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.596"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5530"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.114.32"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01572"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.006"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.545"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8113"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.525"
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.62"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₈123"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.786.19"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4548"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.22"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05292"
$ws.Range("E50").Value = "  +4.12%  "
$ws.Range("E51").Value = "  +0.32%  "
